$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, shifting existing rows 60-135 down to 61-136.
$ws.Rows('60:60').Insert()

# Populate the new row 60 with the data added by the commit (a new weekly
# price observation for Ají - Americana (o), dated 2023-02-09).
$ws.Range('A60').Value = 7
$ws.Range('B60').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C60').Value = 'Ñuble'
$ws.Range('D60').Value = 44966
$ws.Range('E60').Value = 16
$ws.Range('F60').Value = 100112021
$ws.Range('G60').Value = 'Ají'
$ws.Range('H60').Value = 'Americana (o)'
$ws.Range('I60').Value = 'Primera'
$ws.Range('J60').Value = 50
$ws.Range('K60').Value = 8000
$ws.Range('L60').Value = 8000
$ws.Range('M60').Value = 8000
$ws.Range('N60').Value = '$/caja 15 kilos'
$ws.Range('O60').Value = 'Región del Maule'
$ws.Range('P60').Value = 533
$ws.Range('Q60').Value = 15
$ws.Range('R60').Value = 'Hortaliza'
